# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# updates to the Leve profit tables (columns H-N) across all 8 job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 803.13336
$ws.Range("J98").Value = 635
$ws.Range("L98").Value = 635
$ws.Range("N98").Value = -3631

$ws.Range("H111").Value = 3711.6667
$ws.Range("I111").Value = 3711.6667
$ws.Range("K111").Value = 11135.0001
$ws.Range("M111").Value = -8068.000100000001

$ws.Range("H112").Value = 1940.3704
$ws.Range("J112").Value = 2154.0435
$ws.Range("L112").Value = 6462.130500000001
$ws.Range("N112").Value = -8678.130500000001

$ws.Range("H122").Value = 803.13336
$ws.Range("J122").Value = 635
$ws.Range("L122").Value = 1905
$ws.Range("N122").Value = -6805

$ws.Range("H138").Value = 4096.8433
$ws.Range("J138").Value = 4651.2324
$ws.Range("L138").Value = 13953.6972
$ws.Range("N138").Value = -24233.6972

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 322.52942
$ws.Range("I2").Value = 341.46667
$ws.Range("J2").Value = 180.5
$ws.Range("K2").Value = 341.46667
$ws.Range("L2").Value = 180.5
$ws.Range("M2").Value = -228.46667
$ws.Range("N2").Value = -406.5

$ws.Range("H32").Value = 16225.1455
$ws.Range("I32").Value = 6891.1333
$ws.Range("K32").Value = 6891.1333
$ws.Range("M32").Value = -6604.1333

$ws.Range("H45").Value = 1674.3334
$ws.Range("I45").Value = 1674.3334
$ws.Range("K45").Value = 1674.3334
$ws.Range("M45").Value = -1297.3334

$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").ClearContents()

$ws.Range("H116").Value = 322.52942
$ws.Range("I116").Value = 341.46667
$ws.Range("J116").Value = 180.5
$ws.Range("K116").Value = 341.46667
$ws.Range("L116").Value = 180.5
$ws.Range("M116").Value = 1952.53333
$ws.Range("N116").Value = -4768.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 322.52942
$ws.Range("I3").Value = 341.46667
$ws.Range("J3").Value = 180.5
$ws.Range("K3").Value = 341.46667
$ws.Range("L3").Value = 180.5
$ws.Range("M3").Value = -227.46667
$ws.Range("N3").Value = -408.5

$ws.Range("H44").Value = 30000
$ws.Range("J44").Value = 30000
$ws.Range("L44").Value = 30000
$ws.Range("N44").Value = -30994

$ws.Range("H80").Value = 545.3333
$ws.Range("I80").Value = 603.3333
$ws.Range("K80").Value = 603.3333
$ws.Range("M80").Value = 394.6667

$ws.Range("H83").Value = 545.3333
$ws.Range("I83").Value = 603.3333
$ws.Range("K83").Value = 3016.6665
$ws.Range("M83").Value = 1975.3335

$ws.Range("H86").Value = 4481.6
$ws.Range("I86").Value = 4755.2
$ws.Range("K86").Value = 4755.2
$ws.Range("M86").Value = -3632.2

$ws.Range("H89").Value = 4481.6
$ws.Range("I89").Value = 4755.2
$ws.Range("K89").Value = 23776
$ws.Range("M89").Value = -18160

$ws.Range("H105").Value = 4150.037
$ws.Range("I105").Value = 3391.7778
$ws.Range("K105").Value = 3391.7778
$ws.Range("M105").Value = -1644.7778

$ws.Range("H107").Value = 1022.8182
$ws.Range("I107").Value = 1047.3334
$ws.Range("K107").Value = 1047.3334
$ws.Range("M107").Value = 872.6666

$ws.Range("H134").Value = 2187.3794
$ws.Range("I134").Value = 2042
$ws.Range("K134").Value = 6126
$ws.Range("M134").Value = -3591

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3144.5454
$ws.Range("I58").Value = 2301.7144
$ws.Range("J58").Value = 3537.8667
$ws.Range("K58").Value = 2301.7144
$ws.Range("L58").Value = 3537.8667
$ws.Range("M58").Value = -2098.7144
$ws.Range("N58").Value = -3943.8667

$ws.Range("H86").Value = 8313.77
$ws.Range("I86").Value = 5010.857
$ws.Range("K86").Value = 5010.857
$ws.Range("M86").Value = -3887.857

$ws.Range("H89").Value = 8313.77
$ws.Range("I89").Value = 5010.857
$ws.Range("K89").Value = 25054.285
$ws.Range("M89").Value = -19438.285

$ws.Range("H105").Value = 2791.5625
$ws.Range("I105").Value = 820.875
$ws.Range("K105").Value = 820.875
$ws.Range("M105").Value = 926.125

$ws.Range("H134").Value = 3282.074
$ws.Range("I134").Value = 2922.3125
$ws.Range("J134").Value = 3805.3635
$ws.Range("K134").Value = 8766.9375
$ws.Range("L134").Value = 11416.0905
$ws.Range("M134").Value = -6231.9375
$ws.Range("N134").Value = -16486.0905

$ws.Range("H136").Value = 3144.5454
$ws.Range("I136").Value = 2301.7144
$ws.Range("J136").Value = 3537.8667
$ws.Range("K136").Value = 6905.1432
$ws.Range("L136").Value = 10613.6001
$ws.Range("M136").Value = -4355.1432
$ws.Range("N136").Value = -15713.6001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 5000
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 5000
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws.Range("H107").Value = 1318.5454
$ws.Range("I107").Value = 1156.5
$ws.Range("J107").Value = 1354.5555
$ws.Range("K107").Value = 3469.5
$ws.Range("L107").Value = 4063.6665
$ws.Range("M107").Value = -1549.5
$ws.Range("N107").Value = -7903.666499999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3056.5
$ws.Range("I132").Value = 2298.5715
$ws.Range("J132").Value = 4117.6
$ws.Range("K132").Value = 6895.7145
$ws.Range("L132").Value = 12352.8
$ws.Range("M132").Value = -4365.7145
$ws.Range("N132").Value = -17412.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2213.5625
$ws.Range("I16").Value = 2194.3845
$ws.Range("K16").Value = 2194.3845
$ws.Range("M16").Value = -2024.3845

$ws.Range("H46").Value = 3042.1428
$ws.Range("I46").Value = 2441.4285
$ws.Range("J46").Value = 3642.8572
$ws.Range("K46").Value = 2441.4285
$ws.Range("L46").Value = 3642.8572
$ws.Range("M46").Value = -2253.4285
$ws.Range("N46").Value = -4018.8572

$ws.Range("H132").Value = 3777.139
$ws.Range("I132").Value = 3637.4075
$ws.Range("K132").Value = 10912.2225
$ws.Range("M132").Value = -8382.2225

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1307
$ws.Range("I122").Value = 1307
$ws.Range("K122").Value = 3921
$ws.Range("M122").Value = -1471
